$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "29.102.48"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.899.81"
$ws.Range("E3").Value = "  -0.35%  "

Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.18%  "

Set-TextValue $ws.Range("D5") "325.06"
$ws.Range("E5").Value = "  -0.70%  "

Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.27%  "

Set-TextValue $ws.Range("D7") "0.4600"
$ws.Range("E7").Value = "  -0.53%  "

Set-TextValue $ws.Range("D8") "0.3885"
$ws.Range("E8").Value = "  -1.19%  "

Set-TextValue $ws.Range("D9") "0.07867"
$ws.Range("E9").Value = "  -0.93%  "

Set-TextValue $ws.Range("D10") "0.9894"
$ws.Range("E10").Value = "  -1.26%  "

Set-TextValue $ws.Range("D11") "21.88"
$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "1.894.69"
$ws.Range("E12").Value = "  +0.36%  "

Set-TextValue $ws.Range("D13") "5.778"
$ws.Range("E13").Value = "  +0.19%  "

Set-TextValue $ws.Range("D14") "7.050"

Set-TextValue $ws.Range("D15") "0.07011"
$ws.Range("E15").Value = "  +0.93%  "

Set-TextValue $ws.Range("D16") "87.88"
$ws.Range("E16").Value = "  -0.64%  "

Set-TextValue $ws.Range("D17") "1.002"
$ws.Range("E17").Value = "  -0.12%  "

Set-TextValue $ws.Range("D18") "0.000009936"
$ws.Range("E18").Value = "  -1.29%  "

Set-TextValue $ws.Range("D19") "17.03"
$ws.Range("E19").Value = "  -0.68%  "

Set-TextValue $ws.Range("D20") "1.000"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").Value = "29.098.34"
$ws.Range("E21").Value = "  -0.25%  "

Set-TextValue $ws.Range("D22") "5.321"
$ws.Range("E22").Value = "  -0.91%  "

Set-TextValue $ws.Range("D23") "11.10"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "2.099"
$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.090.69"
$ws.Range("E25").Value = "  -1.57%  "

Set-TextValue $ws.Range("D26") "156.18"
$ws.Range("E26").Value = "  -0.17%  "

Set-TextValue $ws.Range("D27") "19.39"
$ws.Range("E27").Value = "  -0.45%  "

Set-TextValue $ws.Range("D28") "5.908"
$ws.Range("E28").Value = "  -2.78%  "

Set-TextValue $ws.Range("D29") "118.50"
$ws.Range("E29").Value = "  -0.37%  "

Set-TextValue $ws.Range("D30") "1.868"
$ws.Range("E30").Value = "  -6.38%  "

Set-TextValue $ws.Range("D31") "0.09316"
$ws.Range("E31").Value = "  -0.81%  "

Set-TextValue $ws.Range("D32") "0.8959"
$ws.Range("E32").Value = "  -3.24%  "

Set-TextValue $ws.Range("D33") "5.238"
$ws.Range("E33").Value = "  -1.76%  "

Set-TextValue $ws.Range("D34") "1.318"
$ws.Range("E34").Value = "  -2.31%  "

Set-TextValue $ws.Range("D35") "3.136"
$ws.Range("E35").Value = "  -4.14%  "

Set-TextValue $ws.Range("D36") "0.05782"
$ws.Range("E36").Value = "  -0.85%  "

Set-TextValue $ws.Range("D37") "1.170"
$ws.Range("E37").Value = "  -3.39%  "

Set-TextValue $ws.Range("D38") "0.02084"
$ws.Range("E38").Value = "  -1.09%  "

Set-TextValue $ws.Range("D39") "1.000"
$ws.Range("E39").Value = "  -0.17%  "

Set-TextValue $ws.Range("D40") "7.661"
$ws.Range("E40").Value = "  -3.54%  "

Set-TextValue $ws.Range("D41") "0.5670"
$ws.Range("E41").Value = "  -1.47%  "

$ws.Range("E42").Value = "  -0.31%  "

Set-TextValue $ws.Range("D43") "9.690"
$ws.Range("E43").Value = "  -2.59%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "11.86"
$ws.Range("E44").Value = "  -1.08%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D45") "2.216"
$ws.Range("E45").Value = "  -0.01%  "

Set-TextValue $ws.Range("D46") "0.5348"
$ws.Range("E46").Value = "  -1.38%  "

Set-TextValue $ws.Range("D47") "0.07007"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("E49").Value = "  -0.65%  "

Set-TextValue $ws.Range("D50") "112.58"
$ws.Range("E50").Value = "  -0.42%  "

Set-TextValue $ws.Range("D51") "1.055"
$ws.Range("E51").Value = "  -2.10%  "
